# Updated cryptos list on Tue Nov 12 10:34:13 UTC 2024 with GitHub Actions
# Refreshes Price (D) and Volume(1h) (E) figures for each coin row, and
# swaps the TRON/Avalanche ordering (rows 13-14) plus replaces the
# Mantle row (51) with ARBITRUM.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.634.64'
$ws.Range("E2").Value = '  +6.89%  '

$ws.Range("D3").Value = '3.332.24'
$ws.Range("E3").Value = '  +4.31%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.97'
$ws.Range("E5").Value = '  +0.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '631.43'
$ws.Range("E6").Value = '  +0.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.401'
$ws.Range("E7").Value = '  +38.80%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.649'
$ws.Range("E8").Value = '  +10.79%  '

$ws.Range("D10").Value = '3.335.83'
$ws.Range("E10").Value = '  +4.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.600'
$ws.Range("E11").Value = '  +1.48%  '

$ws.Range("E12").Value = '  +5.49%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.177'
$ws.Range("E13").Value = '  +7.09%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.43'
$ws.Range("E14").Value = '  +11.97%  '

$ws.Range("D15").Value = '3.947.55'
$ws.Range("E15").Value = '  +4.39%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.35'
$ws.Range("E16").Value = '  +0.50%  '

$ws.Range("D17").Value = '87.505.87'
$ws.Range("E17").Value = '  +6.84%  '

$ws.Range("D18").Value = '3.346.31'
$ws.Range("E18").Value = '  +4.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.43'
$ws.Range("E19").Value = '  +3.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.37'
$ws.Range("E20").Value = '  +5.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '444.79'
$ws.Range("E21").Value = '  +2.25%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.99'
$ws.Range("E22").Value = '  -7.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.46'
$ws.Range("E23").Value = '  +7.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.27'
$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.33'
$ws.Range("E25").Value = '  +1.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.33'
$ws.Range("E26").Value = '  +13.32%  '

$ws.Range("D27").Value = '3.546.27'
$ws.Range("E27").Value = '  +7.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '78.67'
$ws.Range("E28").Value = '  +2.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000136'
$ws.Range("E29").Value = '  +9.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("E31").Value = '  +27.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.06'
$ws.Range("E32").Value = '  +0.14%  '

$ws.Range("E33").Value = '  +0.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '562.44'
$ws.Range("E34").Value = '  -4.42%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("E35").Value = '  -0.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.03'
$ws.Range("E36").Value = '  +2.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.06'
$ws.Range("E37").Value = '  +15.37%  '

$ws.Range("E38").Value = '  -11.13%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.24'
$ws.Range("E39").Value = '  +1.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.420'
$ws.Range("E40").Value = '  +2.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.86'
$ws.Range("E41").Value = '  +5.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.06'
$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.03'
$ws.Range("E44").Value = '  -2.29%  '

$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '156.09'
$ws.Range("E46").Value = '  -2.25%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '182.92'
$ws.Range("E47").Value = '  -2.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.38'
$ws.Range("E48").Value = '  +3.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.78'
$ws.Range("E49").Value = '  +2.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.37'
$ws.Range("E50").Value = '  +3.74%  '

$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.646'
$ws.Range("E51").Value = '  +2.93%  '
